# Update Name of Algo
# Apply updated imputed values in column D for the RandomForest result data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D8"  = -8.5991
    "D10" = -7.840599999999996
    "D12" = -5.923399999999998
    "D18" = -8.6454
    "D37" = -8.123900000000003
    "D55" = -8.949199999999998
    "D68" = -6.933999999999997
    "D77" = -5.674600000000003
    "D78" = -7.8041
    "D81" = -8.059099999999995
    "D82" = -8.246500000000008
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
